$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 32894.25
$ws.Range("J3").Value = 32894.25
$ws.Range("L3").Value = 32894.25
$ws.Range("N3").Value = -33122.25

$ws.Range("H4").Value = 67.63636
$ws.Range("I4").Value = 67.63636
$ws.Range("K4").Value = 67.63636
$ws.Range("M4").Value = 46.36364

$ws.Range("H17").Value = 47791.79
$ws.Range("J17").Value = 50419.11
$ws.Range("L17").Value = 151257.33
$ws.Range("N17").Value = -151593.33

$ws.Range("H28").Value = 1306.8125
$ws.Range("I28").Value = 1060.6
$ws.Range("K28").Value = 1060.6
$ws.Range("M28").Value = -575.5999999999999

$ws.Range("H64").Value = 5734
$ws.Range("I64").Value = 5202
$ws.Range("K64").Value = 5202
$ws.Range("M64").Value = -4954

$ws.Range("H67").Value = 5734
$ws.Range("I67").Value = 5202
$ws.Range("K67").Value = 5202
$ws.Range("M67").Value = -4344

$ws.Range("H100").Value = 2750.6667
$ws.Range("I100").Value = 2424.0908
$ws.Range("J100").Value = 3648.75
$ws.Range("K100").Value = 2424.0908
$ws.Range("L100").Value = 3648.75
$ws.Range("M100").Value = -1883.0908
$ws.Range("N100").Value = -4730.75

$ws.Range("H102").Value = 32894.25
$ws.Range("J102").Value = 32894.25
$ws.Range("L102").Value = 32894.25
$ws.Range("N102").Value = -39384.25

$ws.Range("H107").Value = 1497.5
$ws.Range("I107").Value = 995
$ws.Range("K107").Value = 995
$ws.Range("M107").Value = 925

$ws.Range("H112").Value = 2307.6316
$ws.Range("J112").Value = 2055.5881
$ws.Range("L112").Value = 6166.7643
$ws.Range("N112").Value = -8382.764299999999

$ws.Range("H135").Value = 554.75
$ws.Range("I135").Value = 449.4091
$ws.Range("K135").Value = 4044.6819
$ws.Range("M135").Value = -1509.6819

$ws.Range("H137").Value = 8453.166999999999
$ws.Range("I137").Value = 4606.647
$ws.Range("K137").Value = 13819.941
$ws.Range("M137").Value = -11269.941

$ws.Range("H138").Value = 2727.38
$ws.Range("I138").Value = 1608.9565
$ws.Range("J138").Value = 3680.111
$ws.Range("K138").Value = 4826.8695
$ws.Range("L138").Value = 11040.333
$ws.Range("M138").Value = 313.1305000000002
$ws.Range("N138").Value = -21320.333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 162.45454
$ws.Range("I5").Value = 154.11111
$ws.Range("J5").Value = 200
$ws.Range("K5").Value = 154.11111
$ws.Range("L5").Value = 200
$ws.Range("M5").Value = -42.11111
$ws.Range("N5").Value = -424

$ws.Range("H32").Value = 15042862
$ws.Range("I32").Value = 6974894
$ws.Range("K32").Value = 6974894
$ws.Range("M32").Value = -6974607

$ws.Range("H61").Value = 4618.25
$ws.Range("I61").Value = 4563.7144
$ws.Range("K61").Value = 4563.7144
$ws.Range("M61").Value = -4351.7144

$ws.Range("H63").Value = 1587
$ws.Range("I63").Value = 1587
$ws.Range("K63").Value = 1587
$ws.Range("M63").Value = -901

$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()

$ws.Range("H66").Value = 1587
$ws.Range("I66").Value = 1587
$ws.Range("K66").Value = 7935
$ws.Range("M66").Value = -4503

$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()

$ws.Range("H97").Value = 40889.953
$ws.Range("I97").Value = 1334
$ws.Range("K97").Value = 1334
$ws.Range("M97").Value = -838

$ws.Range("H132").Value = 8575.852000000001
$ws.Range("I132").Value = 8086.5
$ws.Range("J132").Value = 9554.556
$ws.Range("K132").Value = 24259.5
$ws.Range("L132").Value = 28663.668
$ws.Range("M132").Value = -21729.5
$ws.Range("N132").Value = -33723.66800000001

$ws.Range("H136").Value = 4618.25
$ws.Range("I136").Value = 4563.7144
$ws.Range("K136").Value = 13691.1432
$ws.Range("M136").Value = -11141.1432

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 162.45454
$ws.Range("I4").Value = 154.11111
$ws.Range("J4").Value = 200
$ws.Range("K4").Value = 154.11111
$ws.Range("L4").Value = 200
$ws.Range("M4").Value = -39.11111
$ws.Range("N4").Value = -430

$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

$ws.Range("H94").Value = 28253.725
$ws.Range("I94").Value = 1722.5
$ws.Range("K94").Value = 1722.5
$ws.Range("M94").Value = -1271.5

$ws.Range("H99").Value = 44781.5
$ws.Range("I99").Value = 3253.1667
$ws.Range("K99").Value = 3253.1667
$ws.Range("M99").Value = -1755.1667

$ws.Range("H105").Value = 4785942.5
$ws.Range("I105").Value = 6254735.5
$ws.Range("K105").Value = 6254735.5
$ws.Range("M105").Value = -6252988.5

$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()

$ws.Range("H107").Value = 7071.08
$ws.Range("I107").Value = 6507.1055
$ws.Range("J107").Value = 8857
$ws.Range("K107").Value = 6507.1055
$ws.Range("L107").Value = 8857
$ws.Range("M107").Value = -4587.1055
$ws.Range("N107").Value = -12697

$ws.Range("H122").Value = 77998
$ws.Range("J122").Value = 77998
$ws.Range("L122").Value = 77998
$ws.Range("N122").Value = -87798

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4419.115
$ws.Range("I31").Value = 3652.6667
$ws.Range("J31").Value = 5076.0713
$ws.Range("K31").Value = 3652.6667
$ws.Range("L31").Value = 5076.0713
$ws.Range("M31").Value = -3357.6667
$ws.Range("N31").Value = -5666.0713

$ws.Range("H34").Value = 4419.115
$ws.Range("I34").Value = 3652.6667
$ws.Range("J34").Value = 5076.0713
$ws.Range("K34").Value = 3652.6667
$ws.Range("L34").Value = 5076.0713
$ws.Range("M34").Value = -3450.6667
$ws.Range("N34").Value = -5480.0713

$ws.Range("H58").Value = 1384.95
$ws.Range("I58").Value = 1363.1052
$ws.Range("K58").Value = 1363.1052
$ws.Range("M58").Value = -1160.1052

$ws.Range("H59").Value = 17304.46
$ws.Range("J59").Value = 17304.46
$ws.Range("L59").Value = 17304.46
$ws.Range("N59").Value = -19594.46

$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()

$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()

$ws.Range("H105").Value = 1774.4615
$ws.Range("I105").Value = 1829.909
$ws.Range("J105").Value = 1469.5
$ws.Range("K105").Value = 1829.909
$ws.Range("L105").Value = 1469.5
$ws.Range("M105").Value = -82.90900000000011
$ws.Range("N105").Value = -4963.5

$ws.Range("H109").Value = 53249.5
$ws.Range("J109").Value = 53249.5
$ws.Range("L109").Value = 53249.5
$ws.Range("N109").Value = -55329.5

$ws.Range("H132").Value = 4556.1333
$ws.Range("I132").Value = 4560.143
$ws.Range("J132").Value = 4500
$ws.Range("K132").Value = 13680.429
$ws.Range("L132").Value = 13500
$ws.Range("M132").Value = -11150.429
$ws.Range("N132").Value = -18560

$ws.Range("H134").Value = 2909.9546
$ws.Range("I134").Value = 2551
$ws.Range("J134").Value = 6499.5
$ws.Range("K134").Value = 7653
$ws.Range("L134").Value = 19498.5
$ws.Range("M134").Value = -5118
$ws.Range("N134").Value = -24568.5

$ws.Range("H136").Value = 1384.95
$ws.Range("I136").Value = 1363.1052
$ws.Range("K136").Value = 4089.3156
$ws.Range("M136").Value = -1539.3156

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 452.54544
$ws.Range("J23").Value = 457.8
$ws.Range("L23").Value = 1373.4
$ws.Range("N23").Value = -1843.4

$ws.Range("H37").Value = 166750000
$ws.Range("J37").Value = 166750000
$ws.Range("L37").Value = 500250000
$ws.Range("N37").Value = -500250224

$ws.Range("H68").Value = 2451.2856
$ws.Range("J68").Value = 2501.4348
$ws.Range("L68").Value = 7504.3044
$ws.Range("N68").Value = -9126.304400000001

$ws.Range("H71").Value = 2451.2856
$ws.Range("J71").Value = 2501.4348
$ws.Range("L71").Value = 22512.9132
$ws.Range("N71").Value = -30624.9132

$ws.Range("H75").Value = 7092.7144
$ws.Range("I75").Value = 325
$ws.Range("J75").Value = 9799.799999999999
$ws.Range("K75").Value = 975
$ws.Range("L75").Value = 29399.4
$ws.Range("M75").Value = 23
$ws.Range("N75").Value = -31395.4

$ws.Range("H78").Value = 7092.7144
$ws.Range("I78").Value = 325
$ws.Range("J78").Value = 9799.799999999999
$ws.Range("K78").Value = 2925
$ws.Range("L78").Value = 88198.2
$ws.Range("M78").Value = 2067
$ws.Range("N78").Value = -98182.2

$ws.Range("H86").Value = 237.44444
$ws.Range("I86").Value = 237.44444
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 712.33332
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = 473.66668
$ws.Range("N86").ClearContents()

$ws.Range("H89").Value = 237.44444
$ws.Range("I89").Value = 237.44444
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 2136.99996
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = 3791.00004
$ws.Range("N89").ClearContents()

$ws.Range("H107").Value = 1300
$ws.Range("I107").Value = 574.7
$ws.Range("J107").Value = 2105.889
$ws.Range("K107").Value = 1724.1
$ws.Range("L107").Value = 6317.667
$ws.Range("M107").Value = 195.8999999999999
$ws.Range("N107").Value = -10157.667

$ws.Range("H112").Value = 146583.28
$ws.Range("I112").Value = 201210.6
$ws.Range("K112").Value = 603631.8
$ws.Range("M112").Value = -602523.8

$ws.Range("H121").Value = 17650542
$ws.Range("J121").Value = 101905
$ws.Range("L121").Value = 305715
$ws.Range("N121").Value = -308335

$ws.Range("H122").Value = 298
$ws.Range("I122").Value = 150
$ws.Range("J122").Value = 353.5
$ws.Range("K122").Value = 1350
$ws.Range("L122").Value = 3181.5
$ws.Range("M122").Value = 1100
$ws.Range("N122").Value = -8081.5

$ws.Range("H131").Value = 108294.875
$ws.Range("I131").Value = 53115.047
$ws.Range("J131").Value = 213638.19
$ws.Range("K131").Value = 159345.141
$ws.Range("L131").Value = 640914.5700000001
$ws.Range("M131").Value = -154305.141
$ws.Range("N131").Value = -650994.5700000001

$ws.Range("H132").Value = 2001.5625
$ws.Range("I132").Value = 1293.5
$ws.Range("J132").Value = 2102.7144
$ws.Range("K132").Value = 11641.5
$ws.Range("L132").Value = 18924.4296
$ws.Range("M132").Value = -9111.5
$ws.Range("N132").Value = -23984.4296

$ws.Range("H137").Value = 3343.2856
$ws.Range("J137").Value = 5572.6
$ws.Range("L137").Value = 16717.8
$ws.Range("N137").Value = -26917.8

$ws.Range("H140").Value = 3215.5334
$ws.Range("I140").Value = 2941.077
$ws.Range("K140").Value = 8823.231
$ws.Range("M140").Value = -3643.231

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 21065.781
$ws.Range("I97").Value = 716.86365
$ws.Range("J97").Value = 65833.39999999999
$ws.Range("K97").Value = 716.86365
$ws.Range("L97").Value = 65833.39999999999
$ws.Range("M97").Value = -220.86365
$ws.Range("N97").Value = -66825.39999999999

$ws.Range("H126").Value = 43605.875
$ws.Range("I126").Value = 2862.6667
$ws.Range("K126").Value = 8588.000100000001
$ws.Range("M126").Value = -6118.000100000001

$ws.Range("H132").Value = 4304.9023
$ws.Range("I132").Value = 4147.8486
$ws.Range("J132").Value = 4952.75
$ws.Range("K132").Value = 12443.5458
$ws.Range("L132").Value = 14858.25
$ws.Range("M132").Value = -9913.5458
$ws.Range("N132").Value = -19918.25

$ws.Range("H133").Value = 89999.5
$ws.Range("J133").Value = 89999.5
$ws.Range("L133").Value = 89999.5
$ws.Range("N133").Value = -100119.5

$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 24464.264
$ws.Range("I7").Value = 27961.062
$ws.Range("J7").Value = 5814.6665
$ws.Range("K7").Value = 27961.062
$ws.Range("L7").Value = 5814.6665
$ws.Range("M7").Value = -27849.062
$ws.Range("N7").Value = -6038.6665

$ws.Range("H18").Value = 500
$ws.Range("J18").Value = 500
$ws.Range("L18").Value = 500
$ws.Range("N18").Value = -844

$ws.Range("H20").Value = 2092800
$ws.Range("I20").Value = 200000
$ws.Range("J20").Value = 2566000
$ws.Range("K20").Value = 200000
$ws.Range("L20").Value = 2566000
$ws.Range("M20").Value = -199774
$ws.Range("N20").Value = -2566452

$ws.Range("H22").Value = 1229.95
$ws.Range("I22").Value = 1162.4375
$ws.Range("J22").Value = 1500
$ws.Range("K22").Value = 1162.4375
$ws.Range("L22").Value = 1500
$ws.Range("M22").Value = -867.4375
$ws.Range("N22").Value = -2090

$ws.Range("H27").Value = 1229.95
$ws.Range("I27").Value = 1162.4375
$ws.Range("J27").Value = 1500
$ws.Range("K27").Value = 1162.4375
$ws.Range("L27").Value = 1500
$ws.Range("M27").Value = -1055.4375
$ws.Range("N27").Value = -1714

$ws.Range("H46").Value = 26431
$ws.Range("I46").Value = 71670.336
$ws.Range("J46").Value = 1755
$ws.Range("K46").Value = 71670.336
$ws.Range("L46").Value = 1755
$ws.Range("M46").Value = -71482.336
$ws.Range("N46").Value = -2131

$ws.Range("H63").Value = 23085
$ws.Range("J63").Value = 23085
$ws.Range("L63").Value = 23085
$ws.Range("N63").Value = -24583

$ws.Range("H66").Value = 23085
$ws.Range("J66").Value = 23085
$ws.Range("L66").Value = 69255
$ws.Range("N66").Value = -76743

$ws.Range("H68").Value = 36912.75
$ws.Range("I68").Value = 17707.572
$ws.Range("K68").Value = 17707.572
$ws.Range("M68").Value = -16958.572

$ws.Range("H71").Value = 36912.75
$ws.Range("I71").Value = 17707.572
$ws.Range("K71").Value = 88537.86
$ws.Range("M71").Value = -84793.86

$ws.Range("H100").Value = 220000.33
$ws.Range("I100").Value = 106668.664
$ws.Range("J100").Value = 333332
$ws.Range("K100").Value = 106668.664
$ws.Range("L100").Value = 333332
$ws.Range("M100").Value = -106127.664
$ws.Range("N100").Value = -334414

$ws.Range("H126").Value = 24464.264
$ws.Range("I126").Value = 27961.062
$ws.Range("J126").Value = 5814.6665
$ws.Range("K126").Value = 83883.186
$ws.Range("L126").Value = 17443.9995
$ws.Range("M126").Value = -81413.186
$ws.Range("N126").Value = -22383.9995

$ws.Range("H132").Value = 3690.0527
$ws.Range("I132").Value = 3370.9614
$ws.Range("K132").Value = 10112.8842
$ws.Range("M132").Value = -7582.8842

$ws.Range("H136").Value = 3417.0293
$ws.Range("I136").Value = 2468.8076
$ws.Range("K136").Value = 7406.4228
$ws.Range("M136").Value = -4856.4228

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 88878.22
$ws.Range("J62").Value = 143034
$ws.Range("L62").Value = 143034
$ws.Range("N62").Value = -144282

$ws.Range("H65").Value = 88878.22
$ws.Range("J65").Value = 143034
$ws.Range("L65").Value = 715170
$ws.Range("N65").Value = -721410

$ws.Range("H100").Value = 6416.5
$ws.Range("J100").Value = 1999
$ws.Range("L100").Value = 3998
$ws.Range("N100").Value = -5080

$ws.Range("H104").Value = 55185
$ws.Range("J104").Value = 55185
$ws.Range("L104").Value = 55185
$ws.Range("N104").Value = -62173

$ws.Range("H113").Value = 4128.722
$ws.Range("I113").Value = 1722.4166
$ws.Range("J113").Value = 8941.333000000001
$ws.Range("K113").Value = 5167.2498
$ws.Range("L113").Value = 26823.999
$ws.Range("M113").Value = -2997.2498
$ws.Range("N113").Value = -31163.999

$ws.Range("H115").Value = 79999
$ws.Range("J115").Value = 79999
$ws.Range("L115").Value = 79999
$ws.Range("N115").Value = -83133

$ws.Range("H132").Value = 6918.2334
$ws.Range("I132").Value = 7475.8887
$ws.Range("K132").Value = 22427.6661
$ws.Range("M132").Value = -19897.6661
